$d = $word.ActiveDocument

# Update the header date
$d.Content.Find.Execute("2023-09-05 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-06 Wednesday", 2)

# Update each multiplication problem cell (old expression -> new expression).
# Each "old" value is unique within the document, so a plain Find/Replace
# targeting the whole document content is unambiguous for each call.
$d.Content.Find.Execute("15×35=", $true, $false, $false, $false, $false, $true, 1, $false, "18×21=", 2)
$d.Content.Find.Execute("22×93=", $true, $false, $false, $false, $false, $true, 1, $false, "34×82=", 2)
$d.Content.Find.Execute("93×64=", $true, $false, $false, $false, $false, $true, 1, $false, "49×70=", 2)
$d.Content.Find.Execute("60×51=", $true, $false, $false, $false, $false, $true, 1, $false, "62×16=", 2)
$d.Content.Find.Execute("52×73=", $true, $false, $false, $false, $false, $true, 1, $false, "45×36=", 2)
$d.Content.Find.Execute("15×41=", $true, $false, $false, $false, $false, $true, 1, $false, "87×98=", 2)
$d.Content.Find.Execute("61×50=", $true, $false, $false, $false, $false, $true, 1, $false, "63×84=", 2)
$d.Content.Find.Execute("49×95=", $true, $false, $false, $false, $false, $true, 1, $false, "95×18=", 2)
$d.Content.Find.Execute("72×88=", $true, $false, $false, $false, $false, $true, 1, $false, "40×70=", 2)
$d.Content.Find.Execute("86×83=", $true, $false, $false, $false, $false, $true, 1, $false, "32×67=", 2)
$d.Content.Find.Execute("87×26=", $true, $false, $false, $false, $false, $true, 1, $false, "29×68=", 2)
$d.Content.Find.Execute("89×20=", $true, $false, $false, $false, $false, $true, 1, $false, "16×63=", 2)
$d.Content.Find.Execute("29×93=", $true, $false, $false, $false, $false, $true, 1, $false, "37×75=", 2)
$d.Content.Find.Execute("11×19=", $true, $false, $false, $false, $false, $true, 1, $false, "49×19=", 2)
$d.Content.Find.Execute("42×25=", $true, $false, $false, $false, $false, $true, 1, $false, "54×60=", 2)
$d.Content.Find.Execute("95×14=", $true, $false, $false, $false, $false, $true, 1, $false, "36×25=", 2)
$d.Content.Find.Execute("67×37=", $true, $false, $false, $false, $false, $true, 1, $false, "48×44=", 2)
$d.Content.Find.Execute("37×59=", $true, $false, $false, $false, $false, $true, 1, $false, "49×70=", 2)
$d.Content.Find.Execute("97×75=", $true, $false, $false, $false, $false, $true, 1, $false, "55×67=", 2)
$d.Content.Find.Execute("71×97=", $true, $false, $false, $false, $false, $true, 1, $false, "26×87=", 2)
$d.Content.Find.Execute("92×63=", $true, $false, $false, $false, $false, $true, 1, $false, "52×75=", 2)
$d.Content.Find.Execute("42×22=", $true, $false, $false, $false, $false, $true, 1, $false, "83×67=", 2)
$d.Content.Find.Execute("74×91=", $true, $false, $false, $false, $false, $true, 1, $false, "92×76=", 2)
$d.Content.Find.Execute("28×67=", $true, $false, $false, $false, $false, $true, 1, $false, "53×71=", 2)
$d.Content.Find.Execute("47×31=", $true, $false, $false, $false, $false, $true, 1, $false, "27×54=", 2)
